$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column DL. This shifts the existing "Note"
# header (and all following headers/data) one column to the right.
$ws.Columns("DL:DL").Insert()

# Populate the newly inserted column's header.
$ws.Range("DL1").Value2 = "DemonstrationProjectIdentifier"

# Touch the rest of the new column so the (empty) cells are materialized
# in each data row, matching the original sparse layout. Applying a no-op
# border setting forces cell creation without altering any visible style.
$ws.Range("DL2:DL8").Borders.LineStyle = -4142

# Update the Id values (column A) to the new identifiers for each row.
$ws.Range("A2").Value2 = "690148897e79911955eafd42"
$ws.Range("A3").Value2 = "690148897e79911955eafd42"
$ws.Range("A4").Value2 = "690148897e79911955eafd42"
$ws.Range("A5").Value2 = "690148897e79911955eafd42"
$ws.Range("A6").Value2 = "690148897e79911955eafd43"
$ws.Range("A7").Value2 = "690148897e79911955eafd44"
$ws.Range("A8").Value2 = "690148897e79911955eafd44"
